# Scrum meeting attendance roll - apply commit changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-blank row 20 (meeting on 10/10 /1:00) ---
# B20: meeting date/time
$ws.Range("B20").Value = "10/10 /1:00"
# C20: meeting place
$ws.Range("C20").Value = "Google Hangout"
# D20..I20: attendance marks for each team member (A = attend, T = tardy)
$ws.Range("D20").Value = "A"
$ws.Range("E20").Value = "T"
$ws.Range("F20").Value = "A"
$ws.Range("G20").Value = "A"
$ws.Range("H20").Value = "A"
$ws.Range("I20").Value = "A"

# Match the border/format used by the rest of the filled-in data rows
# (copy C19's format onto C20, same as the rest of the table).
$ws.Range("C19").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view: scrolled down with a new active selection ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B21").Select()
